$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header-row labels: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$newHeadersFV2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2410[$i]
}

$newHeadersFV2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

# Column 11 is "diff" and is unchanged; FV2504 headers start at column 12 (L)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2504[$i]
}

# --- 2. Turn the used range A1:U65 into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
